# REPORTGEN-397 : fix extension name in templates
#
# 1) "... metric id from quality model (eg page 4)" - the "eg" was
#    wrapped in spell-check proofErr markers across three runs; collapse
#    that phrase back down to plain text (no spelling-error markers).
# 2) The explanatory paragraph above the "Quality Standards Support"
#    generic-graph sample wrongly called the extension "Standard Quality
#    Rules" - correct it to "Quality Standards Support" (matching the
#    name already used a couple of paragraphs above), which also drags
#    Word's "last edit" (_GoBack) bookmark along with it.

$d = $word.ActiveDocument

# --- (1) " id from quality model (eg page" -----------------------------
# Locate the phrase once and retype exactly that span so the run(s)
# backing "eg" (and its surrounding spell-check markers) collapse back
# into plain text again. A plain re-assignment with identical text is a
# no-op for Word, so clear the span first and then insert the text.
$r = $d.Content
$found = $r.Find.Execute(" id from quality model (eg page", $true, $false, `
    $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $phrase = $d.Range($r.Start, $r.End)
    $phrase.Text = ""
    $phrase.InsertBefore(" id from quality model (eg page")
}

# --- (2) "Standard Quality Rules" -> "Quality Standards Support" -------
# There are two look-alike strings in the document: the heading "...
# sample with Standard Quality Rules" (leave untouched) and the
# explanatory sentence '... extension "Standard Quality Rules" is
# installed ...' (the one that must change). Anchor on the curly-quoted
# phrase, which only the sentence has, so the heading is never touched.
$openQuote = [char]0x201C
$closeQuote = [char]0x201D
$oldPhrase = $openQuote + "Standard Quality Rules" + $closeQuote
$newPhrase = $openQuote + "Quality Standards Support" + $closeQuote

$r2 = $d.Content
$found2 = $r2.Find.Execute($oldPhrase, $true, $false, $false, $false, `
    $false, $true, 1, $false, $newPhrase, 2)

# --- (3) Move the "_GoBack" bookmark ------------------------------------
# Word keeps a single hidden "_GoBack" bookmark at the location of the
# most recent edit. Since the real edit just made is the text fix above,
# drop the old bookmark (currently around the first, untouched
# "Quality Standards Support" occurrence) and re-create it around the
# "Quality Standards Support" text we just typed.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$first = $d.Content
$firstFound = $first.Find.Execute("Quality Standards Support", $true, `
    $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($firstFound) {
    $tail = $d.Range($first.End, $d.Content.End)
    $secondFound = $tail.Find.Execute("Quality Standards Support", $true, `
        $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($secondFound) {
        $d.Bookmarks.Add("_GoBack", $tail) | Out-Null
    }
}
